$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 77, shifting existing rows 77:138 down to 78:139.
$ws.Rows(77).Insert()

# Populate the new row 77 with its final values.
$ws.Range("A77").Value = 11
$ws.Range("B77").Value = 'Vega Monumental Concepción'
$ws.Range("C77").Value = 'Bíobío'
$ws.Range("D77").Value = 44874
$ws.Range("E77").Value = 8
$ws.Range("F77").Value = 'Fruta'
$ws.Range("G77").Value = 100108
$ws.Range("H77").Value = 'Tropicales y subtropicales'
$ws.Range("I77").Value = 100108002
$ws.Range("J77").Value = 'Mango'
$ws.Range("K77").Value = 'Sin especificar'
$ws.Range("L77").Value = 'Primera'
$ws.Range("M77").Value = 200
$ws.Range("N77").Value = 9000
$ws.Range("O77").Value = 9500
$ws.Range("P77").Value = 9250
$ws.Range("Q77").Value = '$/bandeja 4 kilos'
$ws.Range("R77").Value = 'Brasil'
$ws.Range("S77").Value = 2312
$ws.Range("T77").Value = 4
